# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shrink the worker/period detail table.
#    Original table had 19 detail rows (16-34). The refreshed extract
#    only keeps 6 rows. Remove the twelve rows 22-33 (the extra middle
#    entries for OSCAR AGUDELO, ALEXANDER CAMARGO, RAFAEL RANGEL,
#    RODRIGO PEÑATA) then remove the old row 21 (SMITH VASQUEZ DURAN)
#    so the closing "WILSON ARCINIEGAS" row (with its heavier bottom
#    border formatting), originally row 34, slides up into row 21.
#    This also carries the two signature rows up from 39/40 to 26/27.
# ------------------------------------------------------------------
$ws.Range("A22:A33").EntireRow.Delete()
$ws.Range("A21:A21").EntireRow.Delete()

# ------------------------------------------------------------------
# 2) Refresh the summary header figures.
# ------------------------------------------------------------------
$ws.Range("E11").Value2 = 341640      # VALOR MORA
$ws.Range("C13").Value2 = 3           # Cant. Trabajadores
$ws.Range("F13").Value2 = 2           # Cant. Periodos

# ------------------------------------------------------------------
# 3) Rewrite the detail rows with the new worker/period data.
#    Rows 16 and 17 (DANIEL MORELO / JOSEFA GARRIDO, periodo 2507)
#    already hold the correct values after the row deletions above.
# ------------------------------------------------------------------

# Row 18: DANIEL MORELO MORELO NAVARRO, periodo 2507
$ws.Range("C18").Value2 = "73131988"
$ws.Range("D18").Value2 = "DANIEL MORELO MORELO NAVARRO"
$ws.Range("E18").Value2 = "2507"
$ws.Range("F18").Value2 = 56940
$ws.Range("G18").Value2 = 1423500

# Row 19: DANIEL MORELO MORELO NAVARRO, periodo 2508
$ws.Range("C19").Value2 = "73131988"
$ws.Range("D19").Value2 = "DANIEL MORELO MORELO NAVARRO"
$ws.Range("E19").Value2 = "2508"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500

# Row 20: JOSEFA M GARRIDO CASSIANI, periodo 2508
$ws.Range("C20").Value2 = "33335743"
$ws.Range("D20").Value2 = "JOSEFA M GARRIDO CASSIANI"
$ws.Range("E20").Value2 = "2508"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

# Row 21: EULOGIO ANTONIO ARRIETA VILLALBA, periodo 2508 (keeps the
# closing-row style that rode up from the old row 34)
$ws.Range("C21").Value2 = "19890385"
$ws.Range("D21").Value2 = "EULOGIO ANTONIO ARRIETA VILLALBA"
$ws.Range("E21").Value2 = "2508"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500
